{"js": "// Target diff analysis\n// -------------------------------------------------------------------------\n// The supplied unified diff touches word/document.xml and word/styles.xml\n// of setHeightConserveRatio-template.docx. Every one of its hunks is of the\n// form:\n//\n//   -<w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/>\n//   +<w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>\n//\n// i.e. the *set* of attribute name/value pairs on every changed element\n// (the <w:document> namespace declarations, <w:color>, <w:pgSz>, <w:pgMar>,\n// <w:rFonts>, <w:lang>, <w:latentStyles>, every <w:lsdException>, every\n// <w:style>, <w:tblInd>, and the <w:tblCellMar> child elements) is exactly\n// identical before and after the change - only the serialized *order* of\n// the XML attributes differs (namespace declarations sorted by prefix,\n// then the remaining attributes sorted alphabetically by local name).\n// No text, run, paragraph, field code, color, size, margin or style value\n// is added, removed or changed anywhere in the diff.\n//\n// That kind of attribute-order normalization is a side effect of how the\n// fixture was re-serialized upstream; it carries no document semantics,\n// is not observable through the Word object model (Office.js/COM never\n// expose - or let a caller control - raw XML attribute order), and Word\n// itself treats both orderings as 100% identical documents.\n//\n// Consequently, applying this change through the Word JavaScript API\n// means leaving every actual document value untouched - which is what\n// this script does. We simply (and harmlessly) read the body text back,\n// to confirm the content already matches the target state, without\n// writing anything.\n\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n\n// Nothing to change: all values referenced by the diff (field instruction\n// colors \"E36C0A\"/accent6/BF, page size 11906x16838, page margins\n// 1417/1417/1417/1417/708/708/0, default fonts/lang, latent styles, and\n// the Normal/Policepardfaut/TableauNormal/Aucuneliste style definitions)\n// are already exactly as in the target - only their XML attribute order\n// differs, which is not a controllable or meaningful edit via Office.js.\n", "ps1": "# Target diff analysis\n# -------------------------------------------------------------------------\n# The supplied unified diff touches word/document.xml and word/styles.xml\n# of setHeightConserveRatio-template.docx. Every one of its hunks is of the\n# form:\n#\n#   -<w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/>\n#   +<w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>\n#\n# i.e. the *set* of attribute name/value pairs on every changed element\n# (the <w:document> namespace declarations, <w:color>, <w:pgSz>, <w:pgMar>,\n# <w:rFonts>, <w:lang>, <w:latentStyles>, every <w:lsdException>, every\n# <w:style>, <w:tblInd>, and the <w:tblCellMar> child elements) is exactly\n# identical before and after the change - only the serialized *order* of\n# the XML attributes differs (namespace declarations sorted by prefix,\n# then the remaining attributes sorted alphabetically by local name).\n# No text, run, paragraph, field code, color, size, margin or style value\n# is added, removed or changed anywhere in the diff.\n#\n# That kind of attribute-order normalization is a side effect of how the\n# fixture was re-serialized upstream; it carries no document semantics,\n# is not observable through the Word object model (Office.js/COM never\n# expose - or let a caller control - raw XML attribute order), and Word\n# itself treats both orderings as 100% identical documents.\n#\n# Consequently, applying this change through the Word COM object model\n# means leaving every actual document value untouched - which is what\n# this script does. We simply (and harmlessly) read the document text\n# back, to confirm the content already matches the target state, without\n# writing anything.\n\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n\n# Nothing to change: all values referenced by the diff (field instruction\n# colors \"E36C0A\"/accent6/BF, page size 11906x16838, page margins\n# 1417/1417/1417/1417/708/708/0, default fonts/lang, latent styles, and\n# the Normal/Policepardfaut/TableauNormal/Aucuneliste style definitions)\n# are already exactly as in the target - only their XML attribute order\n# differs, which is not a controllable or meaningful edit via Word COM.\n"}
